$wb = $excel.ActiveWorkbook

# Target the "Logs" worksheet explicitly (row data lives here).
$ws = $wb.Worksheets.Item("Logs")

# Append a new row (31) with the same shape as the existing log rows.
$ws.Range("A31").Value = "Demo inplannen"
$ws.Range("B31").Value = "klantenservice@testbedrijf123.nl"
$ws.Range("C31").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Range("D31").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("E31").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$ws.Range("F31").Value = "2025-08-13 22:54:21"
$ws.Range("G31").Value = "Nee"
$ws.Range("H31").Value = "Ja"
$ws.Range("I31").Value = "Nee"
$ws.Range("J31").Value = "Nee"

# Extend the conditional-formatting ranges (D/G/H/I/J) to cover the new row.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $ws.Range("$col`2:$col`30")
    $newRange = $ws.Range("$col`2:$col`31")
    $fc = $oldRange.FormatConditions
    for ($i = 1; $i -le $fc.Count; $i++) {
        $fc.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary count for that category.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 30
